$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format per cell so that values such as "1.00", "59.377.59"
# or "0.0000141" are stored verbatim instead of being auto-converted to
# numbers (which would drop trailing zeros / use scientific notation).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.377.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.522.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.95"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.97"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.575"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.551.84"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.63%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.59"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.60%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.968.36"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.85"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.277.78"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.546.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.27"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.02"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.86"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.98"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.439"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.92"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.73%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.70"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.73"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.69"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.97"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "299.29"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.42%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.831"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.993"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.606"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.78"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0934"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.38"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.82"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0515"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.01%  "
